# feat: add 2022-Q1 data
#
# The workbook has sheets: 2021-Q2, 2021-Q3, 2021-Q4, 总计 (Total).
# This script:
#   1. Inserts a new worksheet "2022-Q1" right before "总计", populated with
#      fund-holding data (same layout as the other quarterly sheets).
#   2. Inserts a new first data row into "总计" summarizing the 2022-Q1
#      quarter, shifting the existing history rows down and renumbering the
#      index column.

$wb = $excel.ActiveWorkbook

# NOTE: inserting a worksheet shifts the position of every sheet that comes
# after it, which invalidates worksheet references captured *before* the
# insertion. So we grab the "total" sheet reference just to use as the
# insertion anchor, perform the Add, and then re-fetch fresh references
# (by name) for everything we touch afterwards.
$totalSheetAnchor = $wb.Worksheets.Item($wb.Worksheets.Count)   # "总计" sheet

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet right before "总计"
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Add($totalSheetAnchor)
$q1Sheet.Name = "2022-Q1"

# Re-fetch fresh references now that the sheet collection has changed.
$templateSheet = $wb.Worksheets.Item("2021-Q4")  # used as a layout template
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# Copy the layout/formatting of the 2021-Q4 sheet (same columns/styles).
# (Column A of row 1 is intentionally left blank, just like the template.)
$templateSheet.Range("B1:H1").Copy($q1Sheet.Range("B1"))
$templateSheet.Range("A2:H3").Copy($q1Sheet.Range("A2"))

# Match the page margins used throughout the rest of the workbook
$q1Sheet.PageSetup.LeftMargin = 54
$q1Sheet.PageSetup.RightMargin = 54
$q1Sheet.PageSetup.TopMargin = 72
$q1Sheet.PageSetup.BottomMargin = 72
$q1Sheet.PageSetup.HeaderMargin = 36
$q1Sheet.PageSetup.FooterMargin = 36

# Force the fund-code / numeric-looking text columns to stay as text
$q1Sheet.Range("B2:G3").NumberFormat = "@"

$q1Sheet.Range("B2").Value = "005104"
$q1Sheet.Range("C2").Value = "富荣福康混合A"
$q1Sheet.Range("D2").Value = "0.08"
$q1Sheet.Range("E2").Value = "87.88"
$q1Sheet.Range("F2").Value = "3.03"
$q1Sheet.Range("G2").Value = "0.0024"
$q1Sheet.Range("H2").Value = 6

$q1Sheet.Range("B3").Value = "005105"
$q1Sheet.Range("C3").Value = "富荣福康混合C"
$q1Sheet.Range("D3").Value = "0.04"
$q1Sheet.Range("E3").Value = "87.88"
$q1Sheet.Range("F3").Value = "3.03"
$q1Sheet.Range("G3").Value = "0.0012"
$q1Sheet.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row at the top of "总计"'s data
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row picks up stray formatting - clear it for the text/number
# columns so the new cells stay unstyled like the other data rows.
$totalSheet.Range("B2:D2").ClearFormats()

# A2 should use the same style as the other index cells (A3, A4, ...)
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0

# Renumber the index column for the shifted-down historical rows
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
